$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F6: 想去人数 (interested count) 4411 -> 4415
    $ws.Range("F6").Value = 4415

    # G8: 最低票价 (lowest price) 0.1 -> "已售罄" (sold out)
    $ws.Range("G8").Value = "已售罄"
}
